$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (columns B..Q), identical for every data row (2..26)
$newValues = @(
    0.9999969141016266,
    0.9990169585194482,
    0.999997262773752,
    0.9999902150793156,
    0.9999933756228668,
    0.00000288055067624043,
    0.0009176260715613002,
    0.00000239901763520588,
    0.00001297651177530373,
    0.000007687764705254807,
    0.00009999963621106398,
    0.00169721851163615,
    0.9999753128130129,
    0.001769472617403162,
    67.51505815013246,
    93.11145047236465
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 2]
    }
}
